$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "QB Website updated 12/13": the Occurrence column (E) previously showed
# "34" for every question in this sheet (row 6 had a stray "34," typo).
# Question 38 now shares the same occurrence, so every populated cell in
# the Occurrence column is updated to read "34, 38".
$ws.Range("E2").Value = "34, 38"
$ws.Range("E3").Value = "34, 38"
$ws.Range("E4").Value = "34, 38"
$ws.Range("E5").Value = "34, 38"
$ws.Range("E6").Value = "34, 38"
$ws.Range("E7").Value = "34, 38"
$ws.Range("E8").Value = "34, 38"
$ws.Range("E9").Value = "34, 38"
$ws.Range("E10").Value = "34, 38"
$ws.Range("E11").Value = "34, 38"
$ws.Range("E12").Value = "34, 38"

# Leave the view scrolled/selected near the edited rows, matching the
# author's saved cursor position.
$ws.Range("E12").Select()
